$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "05/26/2021 01:55:32"
$ws.Range("B2").Value = 24.506

# Update row 3
$ws.Range("A3").Value = "05/26/2021 01:56:34"
$ws.Range("B3").Value = 27.972

# Clear rows 4 through 24 (content only, so the sheet dimension shrinks to A1:B3)
$ws.Range("A4:B24").ClearContents()
